# Insert a new data row before the current row 28 (shifting rows 28-79 down
# to rows 29-80, which matches values that used to be in the row above them),
# then populate the newly inserted row 28 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 28 - this shifts existing rows 28..79 down to 29..80.
$ws.Rows.Item(28).Insert()

# Fill in the new row 28 with the new record (same shape as the other rows in
# the table, i.e. columns A-L, Q and T follow the common template used by all
# other "Macroferia Regional de Talca" / Mango rows).
$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "Macroferia Regional de Talca"
$ws.Range("C28").Value = "Maule"
$ws.Range("D28").Value = 44482
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100108
$ws.Range("H28").Value = "Tropicales y subtropicales"
$ws.Range("I28").Value = 100108002
$ws.Range("J28").Value = "Mango"
$ws.Range("K28").Value = "Sin especificar"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 256
$ws.Range("N28").Value = 7000
$ws.Range("O28").Value = 7000
$ws.Range("P28").Value = 7000
$ws.Range("Q28").Value = "`$/bandeja 4 kilos"
$ws.Range("R28").Value = "Perú"
$ws.Range("S28").Value = 1750
$ws.Range("T28").Value = 4
